# Update "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values on the zh-cn and de-de report sheets (regenerated report timestamps).

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D3").Value = "2016-02-24 07:11:44"
$zhcn.Range("G3").Value = "2016-02-24 07:12:30"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D3").Value = "2016-02-24 07:11:56"
$dede.Range("G3").Value = "2016-02-24 07:12:51"
